# Restore C10 (rule "R30" lower bound) on the "Rules" sheet from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
